$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Append an ORDER BY / LIMIT clause to the Neo4j query stored in B2
# (cell B2 on the "startup" sheet holds the longer Cypher query used by
# the TC05 Bento "EndocrineTherapy-Other" test case).
$ws.Range("B2").Value = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
 WHERE   tp.endocrine_therapy_type IN ["Other"] 
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`
 order By ss.study_subject_id ASC LIMIT 100 
'@

# The view had scrolled so row 2 sat at the top (topLeftCell = B2); restore
# the natural scroll position so row 1 is visible again (topLeftCell = B1).
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2

# Keep the original selection (C2).
$ws.Range("C2").Select() | Out-Null
